$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).ColumnWidth = 29.17

$ws.Range("A1").Value = "id"
$ws.Range("A2").Value = 252681
$ws.Range("A3").Value = 253497
